# Add two new columns (I: "I0", J: "IF") to the worksheet.
# Column I is a constant of 1 for every data row; column J duplicates
# the value already present in column H for that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the last used row based on column A (data rows start at row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Header row: copy the style of the existing header cell (H1) so the new
# header cells I1/J1 look the same as the rest of the header row.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
